$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 6651.0293
$ws.Range("I19").Value = 1747.65
$ws.Range("K19").Value = 1747.65
$ws.Range("M19").Value = -1572.65
$ws.Range("H28").Value = 25641542
$ws.Range("J28").Value = 414.66666
$ws.Range("L28").Value = 414.66666
$ws.Range("N28").Value = -1384.66666
$ws.Range("H33").Value = 687.6875
$ws.Range("I33").Value = 689.5
$ws.Range("K33").Value = 689.5
$ws.Range("M33").Value = -460.5
$ws.Range("H38").Value = 2031.75
$ws.Range("I38").Value = 101.36364
$ws.Range("J38").Value = 6278.6
$ws.Range("K38").Value = 304.09092
$ws.Range("L38").Value = 18835.8
$ws.Range("M38").Value = 67.90908000000002
$ws.Range("N38").Value = -19579.8
$ws.Range("H43").Value = 1695.6875
$ws.Range("I43").Value = 1361.3334
$ws.Range("J43").Value = 2125.5715
$ws.Range("K43").Value = 1361.3334
$ws.Range("L43").Value = 2125.5715
$ws.Range("M43").Value = -1292.3334
$ws.Range("N43").Value = -2263.5715
$ws.Range("H53").Value = 554.8125
$ws.Range("I53").Value = 681.1111
$ws.Range("J53").Value = 392.42856
$ws.Range("K53").Value = 681.1111
$ws.Range("L53").Value = 392.42856
$ws.Range("M53").Value = -44.11109999999996
$ws.Range("N53").Value = -1666.42856
$ws.Range("H116").Value = 6071.8125
$ws.Range("I116").Value = 7095.25
$ws.Range("K116").Value = 7095.25
$ws.Range("M116").Value = -3653.25
$ws.Range("H127").Value = 2548.524
$ws.Range("I127").Value = 2566.15
$ws.Range("K127").Value = 7698.450000000001
$ws.Range("M127").Value = -2738.450000000001
$ws.Range("H138").Value = 2907.8472
$ws.Range("I138").Value = 2713.5789
$ws.Range("J138").Value = 2977.4905
$ws.Range("K138").Value = 8140.736699999999
$ws.Range("L138").Value = 8932.4715
$ws.Range("M138").Value = -3000.736699999999
$ws.Range("N138").Value = -19212.4715

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7921.0254
$ws.Range("I32").Value = 7540.7236
$ws.Range("J32").Value = 17555.334
$ws.Range("K32").Value = 7540.7236
$ws.Range("L32").Value = 17555.334
$ws.Range("M32").Value = -7253.7236
$ws.Range("N32").Value = -18129.334
$ws.Range("H45").Value = 11658.565
$ws.Range("I45").Value = 12088.857
$ws.Range("J45").Value = 10989.223
$ws.Range("K45").Value = 12088.857
$ws.Range("L45").Value = 10989.223
$ws.Range("M45").Value = -11711.857
$ws.Range("N45").Value = -11743.223
$ws.Range("H61").Value = 5349.724
$ws.Range("J61").Value = 4273.625
$ws.Range("L61").Value = 4273.625
$ws.Range("N61").Value = -4697.625
$ws.Range("H132").Value = 3128.625
$ws.Range("I132").Value = 2831.4443
$ws.Range("J132").Value = 5803.25
$ws.Range("K132").Value = 8494.332900000001
$ws.Range("L132").Value = 17409.75
$ws.Range("M132").Value = -5964.332900000001
$ws.Range("N132").Value = -22469.75
$ws.Range("H136").Value = 5349.724
$ws.Range("J136").Value = 4273.625
$ws.Range("L136").Value = 12820.875
$ws.Range("N136").Value = -17920.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 10561.897
$ws.Range("J99").Value = 13936.895
$ws.Range("L99").Value = 13936.895
$ws.Range("N99").Value = -16932.895
$ws.Range("H126").Value = 10561.897
$ws.Range("J126").Value = 13936.895
$ws.Range("L126").Value = 41810.685
$ws.Range("N126").Value = -46750.685
$ws.Range("H138").Value = 143298
$ws.Range("J138").Value = 143298
$ws.Range("L138").Value = 143298
$ws.Range("N138").Value = -153578

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 886.125
$ws.Range("I25").Value = 798.5714
$ws.Range("K25").Value = 2395.7142
$ws.Range("M25").Value = -2226.7142
$ws.Range("H30").Value = 886.125
$ws.Range("I30").Value = 798.5714
$ws.Range("K30").Value = 2395.7142
$ws.Range("M30").Value = -2293.7142
$ws.Range("H34").Value = 698.2857
$ws.Range("I34").Value = 669.6667
$ws.Range("J34").Value = 719.75
$ws.Range("K34").Value = 2009.0001
$ws.Range("L34").Value = 2159.25
$ws.Range("M34").Value = -1925.0001
$ws.Range("N34").Value = -2327.25
$ws.Range("H36").Value = 3348.2
$ws.Range("I36").Value = 580.6667
$ws.Range("K36").Value = 1742.0001
$ws.Range("M36").Value = -1573.0001
$ws.Range("H39").Value = 8950
$ws.Range("I39").Value = 3499.75
$ws.Range("J39").Value = 10161.167
$ws.Range("K39").Value = 10499.25
$ws.Range("L39").Value = 30483.501
$ws.Range("M39").Value = -10205.25
$ws.Range("N39").Value = -31071.501
$ws.Range("H55").Value = 2287.4736
$ws.Range("I55").Value = 2043.6364
$ws.Range("J55").Value = 2622.75
$ws.Range("K55").Value = 6130.9092
$ws.Range("L55").Value = 7868.25
$ws.Range("M55").Value = -5953.9092
$ws.Range("N55").Value = -8222.25
$ws.Range("H88").Value = 9997.5
$ws.Range("I88").Value = 6995
$ws.Range("J88").Value = 13000
$ws.Range("K88").Value = 20985
$ws.Range("L88").Value = 39000
$ws.Range("M88").Value = -20557
$ws.Range("N88").Value = -39856
$ws.Range("H91").Value = 9997.5
$ws.Range("I91").Value = 6995
$ws.Range("J91").Value = 13000
$ws.Range("K91").Value = 20985
$ws.Range("L91").Value = 39000
$ws.Range("M91").Value = -19503
$ws.Range("N91").Value = -41964
$ws.Range("H104").Value = 10231.25
$ws.Range("J104").Value = 13333.333
$ws.Range("L104").Value = 39999.999
$ws.Range("N104").Value = -45241.999
$ws.Range("H122").Value = 645.13336
$ws.Range("J122").Value = 498.16666
$ws.Range("L122").Value = 4483.49994
$ws.Range("N122").Value = -9383.49994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 79351
$ws.Range("I46").Value = 11111
$ws.Range("K46").Value = 11111
$ws.Range("M46").Value = -10955
$ws.Range("H122").Value = 668390.1
$ws.Range("I122").Value = 1112081.6
$ws.Range("J122").Value = 2852.8333
$ws.Range("K122").Value = 3336244.8
$ws.Range("L122").Value = 8558.499899999999
$ws.Range("M122").Value = -3333794.8
$ws.Range("N122").Value = -13458.4999
$ws.Range("H140").Value = 100849.4
$ws.Range("J140").Value = 100849.4
$ws.Range("L140").Value = 100849.4
$ws.Range("N140").Value = -111209.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1947.6207
$ws.Range("I61").Value = 1803.5652
$ws.Range("K61").Value = 1803.5652
$ws.Range("M61").Value = -1601.5652
$ws.Range("H113").Value = 1947.6207
$ws.Range("I113").Value = 1803.5652
$ws.Range("K113").Value = 1803.5652
$ws.Range("M113").Value = 366.4348
$ws.Range("H122").Value = 10989.667
$ws.Range("I122").Value = 14547.667
$ws.Range("K122").Value = 43643.001
$ws.Range("M122").Value = -41193.001
$ws.Range("H132").Value = 57971.145
$ws.Range("I132").Value = 64964.06
$ws.Range("J132").Value = 28251.25
$ws.Range("K132").Value = 194892.18
$ws.Range("L132").Value = 84753.75
$ws.Range("M132").Value = -192362.18
$ws.Range("N132").Value = -89813.75
$ws.Range("H136").Value = 4872151.5
$ws.Range("I136").Value = 10010908
$ws.Range("J136").Value = 3855.4211
$ws.Range("K136").Value = 30032724
$ws.Range("L136").Value = 11566.2633
$ws.Range("M136").Value = -30030174
$ws.Range("N136").Value = -16666.2633

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 15216.833
$ws.Range("I122").Value = 6273.1113
$ws.Range("K122").Value = 18819.3339
$ws.Range("M122").Value = -16369.3339
$ws.Range("H136").Value = 2126.5151
$ws.Range("I136").Value = 1942.28
$ws.Range("J136").Value = 2702.25
$ws.Range("K136").Value = 5826.84
$ws.Range("L136").Value = 8106.75
$ws.Range("M136").Value = -3276.84
$ws.Range("N136").Value = -13206.75
$ws.Range("H139").Value = 68699.5
$ws.Range("J139").Value = 68699.5
$ws.Range("L139").Value = 68699.5
$ws.Range("N139").Value = -78979.5
